# Add the 2022 data column (column S) to the neonatal-mortality table,
# mirroring the formatting already used by the adjacent 2021 column (R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> 2022 value (same order/rows as the existing D:R year columns)
$values = [ordered]@{
    4  = 2022
    5  = 10.071559327675153
    6  = 10.551906067345987
    7  = 9.5619606820956751
    8  = 8.2747510251903922
    9  = 7.6325088339222615
    10 = 8.9652028567087072
    11 = 8.5830821067565175
    12 = 10.275380189066995
    13 = 6.7661984261234096
    14 = 9.0818473806623103
    15 = 9.0186815546489161
    16 = 9.149130832570906
    17 = 8.0270384452893957
    18 = 8.8235294117647065
    19 = 7.2217502124044177
    20 = 3.4213262670647033
    21 = 3.4802022457154114
    22 = 3.3598464070213931
    23 = 12.808072967203572
    24 = 14.988470407378941
    25 = 10.584084672677381
    26 = 7.1442946266854497
    27 = 7.5305623471882646
    28 = 6.7350533623458704
    29 = 16.241806263656226
    30 = 15.614010409340272
    31 = 16.915873735085334
    32 = 22.910065805508165
    33 = 24.889729048519218
    34 = 20.801878879382652
}

foreach ($r in $values.Keys) {
    # Column R is 18, column S is 19 - copy R's formatting onto S first so
    # the new column visually matches the rest of the year columns, then
    # write the 2022 value on top of it.
    $srcCell = $ws.Cells.Item($r, 18)
    $dstCell = $ws.Cells.Item($r, 19)

    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $dstCell.Value = $values[$r]
}

# Clear the clipboard marching ants / mode.
$excel.CutCopyMode = 0

# Reset the selection away from the old "R3" cell so the workbook doesn't
# keep pointing at a cell that's no longer relevant to the edit.
$ws.Range("A1").Select() | Out-Null
